# Adds 5 new data rows (68-72, dates 2024-07-25 .. 2024-07-29) to Sheet1,
# extending the dimension from A1:Z67 to A1:Z72 (matches the target diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column-letter -> column-index map used below (sheet columns A..Z).
$colIndex = @{
    "A" = 1;  "B" = 2;  "C" = 3;  "D" = 4;  "E" = 5;  "F" = 6;  "G" = 7;
    "H" = 8;  "I" = 9;  "J" = 10; "K" = 11; "L" = 12; "M" = 13; "N" = 14;
    "O" = 15; "P" = 16; "Q" = 17; "R" = 18; "S" = 19; "T" = 20; "U" = 21;
    "V" = 22; "W" = 23; "X" = 24; "Y" = 25; "Z" = 26
}

# New row data exactly as it appears in the target worksheet XML. Columns
# H and V are intentionally omitted for every new row (they stay blank,
# same as row 67 right above them).
# Note: scientific-notation literals (e.g. 2.7552E-06) are not accepted
# directly by the PowerShell parser here, so those few values are cast
# from string via [double]"..." instead.
$newRowsData = @(
    @{ Row = 68; Date = 45498; Cells = @{ "B" = 702.298658338; "C" = 220.176192698; "D" = 0; "E" = 0; "F" = 0; "G" = 118.7477283; "I" = 252.0726750185; "J" = 0; "K" = 21.202920191988; "L" = 0; "M" = 0; "N" = 128.41017442016; "O" = 57.958906122; "P" = 0; "Q" = [double]"2.7552E-06"; "R" = 0; "S" = 0; "T" = 0; "U" = 355.8985284435022; "W" = 0; "X" = 0; "Y" = 0; "Z" = 249.925651023048 } }
    @{ Row = 69; Date = 45499; Cells = @{ "B" = 724.7982751876001; "C" = 227.0495050735; "D" = 0; "E" = 0; "F" = 0; "G" = 119.06235135; "I" = 268.5229943598; "J" = 0; "K" = 21.50114173625; "L" = 0; "M" = 0; "N" = 138.15902930048; "O" = 59.005132433; "P" = 0; "Q" = [double]"2.8848E-06"; "R" = 0; "S" = 0; "T" = 0; "U" = 370.7382945468258; "W" = 0; "X" = 0; "Y" = 0; "Z" = 249.36388007919 } }
    @{ Row = 70; Date = 45500; Cells = @{ "B" = 863.2435696150001; "C" = 225.2744945135; "D" = 0; "E" = 0; "F" = 0; "G" = 0; "I" = 269.791815788; "J" = 0; "K" = 0.052158167404; "L" = 0; "M" = 0; "N" = 137.42200247104; "O" = 59.09071703200001; "P" = 0; "Q" = [double]"2.856E-06"; "R" = 0; "S" = 0; "T" = 0; "U" = 409.7566450771163; "W" = 0; "X" = 0; "Y" = 0; "Z" = 249.634362385492 } }
    @{ Row = 71; Date = 45501; Cells = @{ "B" = 867.7364818068002; "C" = 226.740958316; "D" = 0; "E" = 0; "F" = 0; "G" = 0; "I" = 271.236738203; "J" = 0; "K" = 0.053829917064; "L" = 0; "M" = 0; "N" = 135.04341588512; "O" = 59.324397653; "P" = 0; "Q" = [double]"2.82E-06"; "R" = 0; "S" = 0; "T" = 0; "U" = 405.7908455150211; "W" = 0; "X" = 0; "Y" = 0; "Z" = 239.293615752254 } }
    @{ Row = 72; Date = 45502; Cells = @{ "B" = 849.1078949759001; "C" = 230.034434941; "D" = 0; "E" = 0; "F" = 0; "G" = 0; "I" = 267.8260438353; "J" = 0; "K" = 0.054533175125; "L" = 0; "M" = 0; "N" = 136.3499634464; "O" = 58.349035061; "P" = 0; "Q" = [double]"2.8296E-06"; "R" = 0; "S" = 0; "T" = 0; "U" = 393.8934468287359; "W" = 0; "X" = 0; "Y" = 0; "Z" = 243.891814959388 } }
)

# Row 67 column A already carries the "date" cell style (bordered, bold,
# centered, custom YYYY-MM-DD HH:MM:SS number format). Copy that style
# across to the new column-A cells so the new rows match it exactly
# instead of minting a brand-new style index.
$ws.Range("A67").Copy() | Out-Null
$ws.Range("A68:A72").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

foreach ($rowData in $newRowsData) {
    $r = $rowData.Row

    # Column A: date serial value (style already copied above).
    $ws.Cells.Item($r, $colIndex["A"]).Value = $rowData.Date

    foreach ($col in $rowData.Cells.Keys) {
        $c = $colIndex[$col]
        $ws.Cells.Item($r, $c).Value = $rowData.Cells[$col]
    }
}
